$wb = $excel.ActiveWorkbook

# --- DLC_List sheet: Duration 600 -> 700, Seed ranges 50 seeds -> 99 seeds ---
$dlc = $wb.Worksheets.Item("DLC_List")

$seedStarts = @(401, 601, 801, 1001, 1201, 1401, 1601, 1801, 2001, 2201, 2401)

for ($i = 0; $i -lt $seedStarts.Count; $i++) {
    $r = $i + 2
    $start = $seedStarts[$i]
    $end = $start + 98
    $dlc.Cells.Item($r, 8).Value = "700"
    $dlc.Cells.Item($r, 10).Value = "[" + $start + ":1:" + $end + "]"
}

# --- Sheet1: widen column P (16) to fit, and mark selection ---
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Columns.Item(16).ColumnWidth = 9.75

# --- Update selections on each sheet (config, Sheet1, then back to DLC_List last) ---
$config = $wb.Worksheets.Item("config")
$config.Activate()
$config.Range("B29").Select()

$sheet1.Activate()
$sheet1.Range("J14").Select()

$dlc.Activate()
$dlc.Range("J12").Select()
